# Applies the cryptos list update described in the commit:
# "Updated cryptos list on Sun Oct  1 03:08:55 UTC 2023 with GitHub Actions"
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.062.68"
$ws.Range("E2").Value = "  +0.53%  "

$ws.Range("D3").Value = "1.675.67"
$ws.Range("E3").Value = "  +0.24%  "

$ws.Range("E4").Value = "  +0.14%  "

$ws.Range("D5").Value = "'215.18"
$ws.Range("E5").Value = "  +0.18%  "

$ws.Range("D6").Value = "'0.516"
$ws.Range("E6").Value = "  -0.12%  "

$ws.Range("E7").Value = "  +0.12%  "

$ws.Range("E8").Value = "  +2.18%  "

$ws.Range("B9").Value = "Dogecoin"
$ws.Range("C9").Value = "https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge"
$ws.Range("D9").Value = "'0.0622"
$ws.Range("E9").Value = "  +0.18%  "

$ws.Range("B10").Value = "Solana"
$ws.Range("C10").Value = "https://coinranking.com/coin/zNZHO_Sjf+solana-sol"
$ws.Range("D10").Value = "'21.23"
$ws.Range("E10").Value = "  +4.62%  "

$ws.Range("D11").Value = "'0.0882"
$ws.Range("E11").Value = "  -0.79%  "

$ws.Range("D12").Value = "1.912.01"
$ws.Range("E12").Value = "  +0.27%  "

$ws.Range("D13").Value = "1.673.50"
$ws.Range("E13").Value = "  -0.34%  "

$ws.Range("E14").Value = "  +0.92%  "

$ws.Range("D15").Value = "'0.535"
$ws.Range("E15").Value = "  +1.55%  "

$ws.Range("D16").Value = "'66.01"
$ws.Range("E16").Value = "  +0.75%  "

$ws.Range("D17").Value = "27.049.71"
$ws.Range("E17").Value = "  +0.49%  "

$ws.Range("D18").Value = "'236.88"
$ws.Range("E18").Value = "  +1.43%  "

$ws.Range("D19").Value = "'8.14"
$ws.Range("E19").Value = "  +1.30%  "

$ws.Range("D20").Value = "0.0₃0740"
$ws.Range("E20").Value = "  +0.94%  "

$ws.Range("E21").Value = "  +0.11%  "

$ws.Range("E22").Value = "  +0.78%  "

$ws.Range("D23").Value = "'9.32"
$ws.Range("E23").Value = "  +1.92%  "

$ws.Range("E24").Value = "  -2.07%  "

$ws.Range("D25").Value = "'146.02"
$ws.Range("E25").Value = "  -0.09%  "

$ws.Range("E26").Value = "  +1.32%  "

$ws.Range("D27").Value = "'16.34"
$ws.Range("E27").Value = "  +2.38%  "

$ws.Range("D28").Value = "'0.113"

$ws.Range("E29").Value = "  +0.18%  "

$ws.Range("E30").Value = "  -0.11%  "

$ws.Range("E31").Value = "  -0.06%  "

$ws.Range("D32").Value = "'3.36"
$ws.Range("E32").Value = "  +0.84%  "

$ws.Range("D33").Value = "1.536.59"
$ws.Range("E33").Value = "  +5.31%  "

$ws.Range("E34").Value = "  +1.51%  "

$ws.Range("D35").Value = "'1.70"
$ws.Range("E35").Value = "  +3.66%  "

$ws.Range("D37").Value = "'0.595"
$ws.Range("E37").Value = "  +2.27%  "

$ws.Range("D38").Value = "'0.917"
$ws.Range("E38").Value = "  +1.75%  "

$ws.Range("E39").Value = "  +2.05%  "

$ws.Range("E40").Value = "  +2.97%  "

$ws.Range("E41").Value = "  +0.11%  "

$ws.Range("E42").Value = "  +1.91%  "

$ws.Range("D43").Value = "'5.54"
$ws.Range("E43").Value = "  -3.48%  "

$ws.Range("E44").Value = "  -1.79%  "

$ws.Range("D45").Value = "1.819.33"
$ws.Range("E45").Value = "  +0.71%  "

$ws.Range("E46").Value = "  +0.29%  "

$ws.Range("D47").Value = "'90.85"
$ws.Range("E47").Value = "  +0.11%  "

$ws.Range("E48").Value = "  +1.17%  "

$ws.Range("E49").Value = "  +2.29%  "

$ws.Range("D50").Value = "'8.07"
$ws.Range("E50").Value = "  +5.83%  "

$ws.Range("E51").Value = "  +0.64%  "
